$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 45209
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("B19").Value = "21:13"
$ws.Range("C19").Value = 1829
$ws.Range("D19").Value = "amazon"
$ws.Range("E19").Value = "preto"
